$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Regenerated sval data (filtered save games) - updated B:G for rows 2-6

$ws.Range("B2").Value = 3.230985683306322
$ws.Range("C2").Value = 1.667794583268128
$ws.Range("D2").Value = 0.8054896365839992
$ws.Range("E2").Value = 0.496779210170732
$ws.Range("G2").Value = 6.201049113329182

$ws.Range("B3").Value = 3.230985683306322
$ws.Range("C3").Value = 1.667794583268128
$ws.Range("D3").Value = 3.900430680208489
$ws.Range("E3").Value = 8.660232485948974
$ws.Range("G3").Value = 17.45944343273191

$ws.Range("B4").Value = 0.3048080303191223
$ws.Range("C4").Value = 0.3127903958511391
$ws.Range("D4").Value = 0.8054896365839992
$ws.Range("E4").Value = 616238.5361209477
$ws.Range("G4").Value = 616239.9592090105

$ws.Range("B5").Value = 0.3048080303191223
$ws.Range("C5").Value = 1.667794583268128
$ws.Range("D5").Value = 0.8054896365839992
$ws.Range("E5").Value = 8.660232485948974
$ws.Range("G5").Value = 11.43832473612022

$ws.Range("B6").Value = 3.230985683306322
$ws.Range("C6").Value = 1.667794583268128
$ws.Range("D6").Value = 26.21740644021617
$ws.Range("E6").Value = 8.660232485948974
$ws.Range("G6").Value = 39.7764191927396
